$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update computed metric columns (M..T) for rows 2-5 with new TPM-derived values
$ws.Range("M2").Value = 30.801072
$ws.Range("N2").Value = 61.602144
$ws.Range("O2").Value = 0.5373480691764108
$ws.Range("P2").Value = 0.485871843331092
$ws.Range("Q2").Value = 17.696386304736
$ws.Range("R2").Value = 106.178317828416
$ws.Range("S2").Value = 0.5373480691764108
$ws.Range("T2").Value = 0.485871843331092

$ws.Range("O3").Value = 0.2096145064786482
$ws.Range("P3").Value = 0.2843011610923331
$ws.Range("S3").Value = 0.2096145064786482
$ws.Range("T3").Value = 0.2843011610923331

$ws.Range("M4").Value = 0.1305583333333333
$ws.Range("N4").Value = 0.391675
$ws.Range("O4").Value = 0.002277689176907768
$ws.Range("P4").Value = 0.003089240761436898
$ws.Range("Q4").Value = 0.07501072371666667
$ws.Range("R4").Value = 0.67509651345
$ws.Range("S4").Value = 0.002277689176907768
$ws.Range("T4").Value = 0.003089240761436898

$ws.Range("M5").Value = 14.373679
$ws.Range("N5").Value = 28.747358
$ws.Range("O5").Value = 0.2507597351680332
$ws.Range("P5").Value = 0.2267377548151379
$ws.Range("Q5").Value = 8.258224785302
$ws.Range("R5").Value = 49.549348711812
$ws.Range("S5").Value = 0.2507597351680332
$ws.Range("T5").Value = 0.2267377548151379

# Remove rows for Neutrophils / Resolving-Mac target clusters (rows 6 and 7)
$ws.Rows("6:7").Delete() | Out-Null
